# Update cached market-price / profit figures across the leve-crafting sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect the latest pull from the
# scheduled market-data runner.
$wb = $excel.ActiveWorkbook

# ===================== Sheet ALC =====================
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 1077.6
$ws.Range("I11").Value = 1077.6
$ws.Range("K11").Value = 1077.6
$ws.Range("M11").Value = -937.5999999999999
# Row 64
$ws.Range("H64").Value = 4689.9
$ws.Range("I64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("M64").Value = -4752
# Row 67
$ws.Range("H67").Value = 4689.9
$ws.Range("I67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("M67").Value = -4142
# Row 69
$ws.Range("H69").Value = 3000
$ws.Range("I69").Value = 3000
$ws.Range("K69").Value = 9000
$ws.Range("M69").Value = -8126
# Row 72
$ws.Range("H72").Value = 3000
$ws.Range("I72").Value = 3000
$ws.Range("K72").Value = 27000
$ws.Range("M72").Value = -22632
# Row 94
$ws.Range("H94").Value = 1483.3334
$ws.Range("I94").Value = 1483.3334
$ws.Range("K94").Value = 1483.3334
$ws.Range("M94").Value = -1032.3334
# Row 101
$ws.Range("H101").Value = 457.7
$ws.Range("J101").Value = 382.5
$ws.Range("L101").Value = 1147.5
$ws.Range("N101").Value = -4391.5
# Row 103
$ws.Range("H103").Value = 950.3333
$ws.Range("I103").Value = 750
$ws.Range("J103").Value = 1050.5
$ws.Range("K103").Value = 2250
$ws.Range("L103").Value = 3151.5
$ws.Range("M103").Value = -1664
$ws.Range("N103").Value = -4323.5
# Row 137
$ws.Range("H137").Value = 1579
$ws.Range("I137").Value = 1226.5
$ws.Range("K137").Value = 3679.5
$ws.Range("M137").Value = -1129.5

# ===================== Sheet ARM =====================
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 900
$ws.Range("I45").Value = 900
$ws.Range("K45").Value = 900
$ws.Range("M45").Value = -523
# Row 74
$ws.Range("H74").Value = 1984
$ws.Range("I74").Value = 1984
$ws.Range("K74").Value = 1984
$ws.Range("M74").Value = -1110
# Row 77
$ws.Range("H77").Value = 1984
$ws.Range("I77").Value = 1984
$ws.Range("K77").Value = 9920
$ws.Range("M77").Value = -5552
# Row 102
$ws.Range("H102").Value = 2799.5
$ws.Range("I102").Value = 2799.5
$ws.Range("K102").Value = 2799.5
$ws.Range("M102").Value = -1177.5

# ===================== Sheet BSM =====================
$ws = $wb.Worksheets.Item("BSM")
# Row 53
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
# Row 82
$ws.Range("H82").Value = 21250.54
# Row 85
$ws.Range("H85").Value = 21250.54
# Row 86
$ws.Range("H86").Value = 2524.2778
$ws.Range("I86").Value = 2571.4167
$ws.Range("J86").Value = 2430
$ws.Range("K86").Value = 2571.4167
$ws.Range("L86").Value = 2430
$ws.Range("M86").Value = -1448.4167
$ws.Range("N86").Value = -4676
# Row 89
$ws.Range("H89").Value = 2524.2778
$ws.Range("I89").Value = 2571.4167
$ws.Range("J89").Value = 2430
$ws.Range("K89").Value = 12857.0835
$ws.Range("L89").Value = 12150
$ws.Range("M89").Value = -7241.083500000001
$ws.Range("N89").Value = -23382
# Row 134
$ws.Range("H134").Value = 5166.2583
$ws.Range("I134").Value = 5115.1724
$ws.Range("K134").Value = 15345.5172
$ws.Range("M134").Value = -12810.5172

# ===================== Sheet CRP =====================
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2390.5293
$ws.Range("I58").Value = 1410.1111
$ws.Range("K58").Value = 1410.1111
$ws.Range("M58").Value = -1207.1111
# Row 122
$ws.Range("H122").Value = 1871
$ws.Range("I122").Value = 1862.4
$ws.Range("K122").Value = 5587.200000000001
$ws.Range("M122").Value = -3137.200000000001
# Row 132
$ws.Range("H132").Value = 3923.75
$ws.Range("J132").Value = 5000
$ws.Range("L132").Value = 15000
$ws.Range("N132").Value = -20060
# Row 136
$ws.Range("H136").Value = 2390.5293
$ws.Range("I136").Value = 1410.1111
$ws.Range("K136").Value = 4230.3333
$ws.Range("M136").Value = -1680.3333

# ===================== Sheet CUL =====================
$ws = $wb.Worksheets.Item("CUL")
# Row 117
$ws.Range("H117").Value = 567
$ws.Range("J117").Value = 466
$ws.Range("L117").Value = 1398
$ws.Range("N117").Value = -8282

# ===================== Sheet GSM =====================
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 199.63637
$ws.Range("I2").Value = 293
$ws.Range("J2").Value = 87.59999999999999
$ws.Range("K2").Value = 293
$ws.Range("L2").Value = 87.59999999999999
$ws.Range("M2").Value = -180
$ws.Range("N2").Value = -313.6
# Row 53
$ws.Range("H53").Value = 30001
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 30001
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 30001
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -31263
# Row 97
$ws.Range("H97").Value = 596.2857
$ws.Range("I97").Value = 596.2857
$ws.Range("K97").Value = 596.2857
$ws.Range("M97").Value = -100.2857

# ===================== Sheet LTW =====================
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 14725.25
$ws.Range("I7").Value = 12345.3
$ws.Range("K7").Value = 12345.3
$ws.Range("M7").Value = -12233.3
# Row 16
$ws.Range("H16").Value = 726.8
$ws.Range("I16").Value = 561.6667
$ws.Range("K16").Value = 561.6667
$ws.Range("M16").Value = -391.6667
# Row 46
$ws.Range("H46").Value = 3967.4119
$ws.Range("I46").Value = 3271.889
$ws.Range("J46").Value = 4749.875
$ws.Range("K46").Value = 3271.889
$ws.Range("L46").Value = 4749.875
$ws.Range("M46").Value = -3083.889
$ws.Range("N46").Value = -5125.875
# Row 57
$ws.Range("H57").Value = 9027.333000000001
$ws.Range("I57").Value = 1041
$ws.Range("K57").Value = 1041
$ws.Range("M57").Value = -475
# Row 100
$ws.Range("H100").Value = 1988.7
$ws.Range("I100").Value = 2099.2222
$ws.Range("J100").Value = 994
$ws.Range("K100").Value = 2099.2222
$ws.Range("L100").Value = 994
$ws.Range("M100").Value = -1558.2222
$ws.Range("N100").Value = -2076
# Row 126
$ws.Range("H126").Value = 14725.25
$ws.Range("I126").Value = 12345.3
$ws.Range("K126").Value = 37035.89999999999
$ws.Range("M126").Value = -34565.89999999999
# Row 132
$ws.Range("H132").Value = 4971
$ws.Range("I132").Value = 4995
$ws.Range("J132").Value = 4899
$ws.Range("K132").Value = 14985
$ws.Range("L132").Value = 14697
$ws.Range("M132").Value = -12455
$ws.Range("N132").Value = -19757
# Row 136
$ws.Range("H136").Value = 3300
$ws.Range("I136").Value = 2450
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 7350
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -4800
$ws.Range("N136").Value = -20100

# ===================== Sheet WVR =====================
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 20000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 20000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 20000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -21040
# Row 113
$ws.Range("H113").Value = 459.6
$ws.Range("J113").Value = 550
$ws.Range("L113").Value = 1650
$ws.Range("N113").Value = -5990
# Row 126
$ws.Range("H126").Value = 22521.53
$ws.Range("I126").Value = 16489
$ws.Range("K126").Value = 49467
$ws.Range("M126").Value = -46997
# Row 132
$ws.Range("H132").Value = 3249.6667
$ws.Range("J132").Value = 3499.5
$ws.Range("L132").Value = 10498.5
$ws.Range("N132").Value = -15558.5
# Row 133
$ws.Range("H133").Value = 150000
$ws.Range("J133").Value = 150000
$ws.Range("L133").Value = 150000
$ws.Range("N133").Value = -160120

